# Update slides 2-7 (1-based Slides collection index): replace the generic
# third-person pronouns "他"/"你" with the reverent, deity-specific forms
# "祂"/"袮" used for God throughout the Psalm text. The substitutions are
# applied character-by-character via TextRange.Characters(start,length) so
# that the run is split at each edit point, mirroring how PowerPoint itself
# records an in-place correction as separate runs.

function Set-Segments($paraRange, [string[]]$segments) {
    # Re-types a paragraph's text as a sequence of segments, each becoming
    # its own run (but inheriting the formatting already on that span).
    $pos = 1
    foreach ($seg in $segments) {
        $len = $seg.Length
        $chars = $paraRange.Characters($pos, $len)
        $chars.Text = $seg
        $pos += $len
    }
}

$p = $ppt.ActivePresentation

# Slide 2 ("耶和華作王！...你的寶座從太初立定；你從亙古就有。")
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
Set-Segments ($tr2.Paragraphs(1,1)) @("耶和華作王", "！", "祂", "以", "威嚴為衣穿上；耶和華以能力為衣，以能力束腰，世界就堅定，不得動搖。")
Set-Segments ($tr2.Paragraphs(2,1)) @("袮", "的", "寶座從太初立定", "；", "袮", "從", "亙古就有。")

# Slide 3 ("耶和華啊，你的法度最的確；你的殿永稱為聖，是合宜的。")
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
Set-Segments ($tr3.Paragraphs(2,1)) @("耶和華啊", "，", "袮", "的", "法度最的確", "；", "袮", "的", "殿永稱為聖，是合宜的。")

# Slide 4 ("他把地建立在海上..." / "誰能站在他的聖所？")
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(1).TextFrame.TextRange
Set-Segments ($tr4.Paragraphs(2,1)) @("祂", "把", "地建立在海上，安定在大水之上。")
Set-Segments ($tr4.Paragraphs(3,1)) @("誰能登耶和華的山？誰能站", "在", "祂", "的", "聖所？")

# Slide 5 ("密雲和幽暗在他的四圍；公義和公平是他寶座的根基。")
$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(1).TextFrame.TextRange
Set-Segments ($tr5.Paragraphs(3,1)) @("密雲和幽暗", "在", "祂", "的", "四圍；公義和公平", "是", "祂", "寶", "座的根基。")

# Slide 6 ("有烈火在他前頭行，燒滅他四圍的敵人。他的閃電光照世界...")
$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(1).TextFrame.TextRange
Set-Segments ($tr6.Paragraphs(1,1)) @("有烈火", "在", "祂", "前", "頭行，燒", "滅", "祂", "四", "圍的敵人", "。", "祂", "的", "閃電光照世界，大地看見便震動。")

# Slide 7 ("諸天表明他的公義；萬民看見他的榮耀。")
$s7 = $p.Slides.Item(7)
$tr7 = $s7.Shapes.Item(1).TextFrame.TextRange
Set-Segments ($tr7.Paragraphs(1,1)) @("諸天表", "明", "祂", "的", "公義；萬民看", "見", "祂", "的", "榮耀。")
